$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_suite")

$data = @(
    @("TCID", "Runmode"),
    @("Login_BillingPortal", "Y"),
    @("View_Dashboard", "Y"),
    @("View_Reports_Daily", "Y"),
    @("Filter_Daily_Reports", "Y"),
    @("Export_Daily_Reports", "Y"),
    @("View_Reports_Monthly", "Y"),
    @("Filter_Monthly_Reports", "Y"),
    @("View_Reports_Yearly", "Y"),
    @("Filter_Yearly_Reports", "Y"),
    @("View_Reports_Custom", "Y"),
    @("Filter_Custom_Reports", "Y"),
    @("Search_Reports", "Y"),
    @("Download_Reports", "Y"),
    @("View_Map", "Y"),
    @("View_Profile", "Y"),
    @("Logout_LogoutOption", "Y"),
    @("ResetPassword_viaForgotPassword", "Y"),
    @("Base", "Y")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("B19").Select()
